# Updates on test data file; igp & vrf variations.
#
# 1. Devices sheet: fill column F ("igp RID") for every device row with the
#    same value as its MgmtIP (column D).
# 2. VRF sheet: populate with rd / import-export route-target data for the
#    four PE routers.
# 3. New "LIB" sheet appended at the end of the workbook (becomes the active
#    tab), with a small device/interface/protocol lookup header row.
# 4. Selections/active-cell bookkeeping updated to match where the author
#    ended up after each edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Devices!F2:F20 = Devices!D2:D20  (igp RID mirrors MgmtIP for test data)
# ---------------------------------------------------------------------
$devices = $wb.Worksheets.Item("Devices")

$devices.Cells.Item(2,6).Value  = "10.224.0.1"
$devices.Cells.Item(3,6).Value  = "10.224.0.2"
$devices.Cells.Item(4,6).Value  = "10.224.0.1"
$devices.Cells.Item(5,6).Value  = "10.224.0.2"
$devices.Cells.Item(6,6).Value  = "10.224.0.3"
$devices.Cells.Item(7,6).Value  = "10.224.0.4"
$devices.Cells.Item(8,6).Value  = "10.224.0.6"
$devices.Cells.Item(9,6).Value  = "10.224.0.7"
$devices.Cells.Item(10,6).Value = "10.224.0.8"
$devices.Cells.Item(11,6).Value = "10.224.0.9"
$devices.Cells.Item(12,6).Value = "10.200.1.4"
$devices.Cells.Item(13,6).Value = "10.200.1.5"
$devices.Cells.Item(14,6).Value = "10.200.10.4"
$devices.Cells.Item(15,6).Value = "10.224.1.1"
$devices.Cells.Item(16,6).Value = "10.224.1.2"
$devices.Cells.Item(17,6).Value = "10.224.1.3"
$devices.Cells.Item(18,6).Value = "10.224.1.4"
$devices.Cells.Item(19,6).Value = "10.224.2.1"
$devices.Cells.Item(20,6).Value = "10.224.2.2"

$devices.Range("H28").Select()

# ---------------------------------------------------------------------
# 2. VRF sheet content
# ---------------------------------------------------------------------
$vrf = $wb.Worksheets.Item("VRF")

# column C (rd) first: header then the four route-distinguisher values
$vrf.Cells.Item(1,3).Value = "rd"
$vrf.Cells.Item(2,3).Value = "10.224.0.1:100"
$vrf.Cells.Item(3,3).Value = "10.224.0.2:100"
$vrf.Cells.Item(4,3).Value = "10.224.0.3:100"
$vrf.Cells.Item(5,3).Value = "10.224.0.4:100"

# sample route-target value, then the route-target header row
$vrf.Cells.Item(2,4).Value = "65400:100"
$vrf.Cells.Item(1,4).Value = "ipv4 import RT"
$vrf.Cells.Item(1,5).Value = "ipv4 export RT"
$vrf.Cells.Item(1,6).Value = "ipv6 import RT"
$vrf.Cells.Item(1,7).Value = "ipv6 export RT"

# vrf name column
$vrf.Cells.Item(2,2).Value = "Test"

# device column
$vrf.Cells.Item(1,1).Value = "device"
$vrf.Cells.Item(1,2).Value = "vrf"
$vrf.Cells.Item(2,1).Value = "scd-pe01"
$vrf.Cells.Item(3,1).Value = "scd-pe02"
$vrf.Cells.Item(4,1).Value = "mtr-pe01"
$vrf.Cells.Item(5,1).Value = "mtr-pe02"

# remaining route-target cells (all reuse existing strings already interned above)
$vrf.Cells.Item(3,2).Value = "Test"
$vrf.Cells.Item(4,2).Value = "Test"
$vrf.Cells.Item(5,2).Value = "Test"

for ($r = 2; $r -le 5; $r++) {
  for ($c = 5; $c -le 7; $c++) {
    $vrf.Cells.Item($r,$c).Value = "65400:100"
  }
  $vrf.Cells.Item($r,4).Value = "65400:100"
}

$vrf.Columns.Item(3).ColumnWidth = 15
$vrf.Columns.Item(4).ColumnWidth = 13
$vrf.Columns.Item(5).ColumnWidth = 12.83
$vrf.Columns.Item(6).ColumnWidth = 13
$vrf.Columns.Item(7).ColumnWidth = 12.83

$vrf.Range("C15").Select()

# ---------------------------------------------------------------------
# 3. New LIB sheet, appended after VRF (becomes the active tab)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lib = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$lib.Name = "LIB"

$lib.Cells.Item(1,1).Value = "Device"
$lib.Cells.Item(1,2).Value = "Interface"
$lib.Cells.Item(1,5).Value = "area"
$lib.Cells.Item(1,3).Value = "label proto"
$lib.Cells.Item(1,4).Value = "igp proto"

$lib.Columns.Item(3).ColumnWidth = 12.83

$lib.Range("E4").Select()
$lib.Activate()
